$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-11-25 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-26 Tuesday", 2)

# Update the equation cells in the practice table.
# Each cell is addressed by (row, column) to avoid ambiguity between
# equations that share text (e.g. "362×9=" appears both as an old value
# and as a new value in different cells).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "991×5="
$t.Cell(1, 2).Range.Text  = "960×8="
$t.Cell(1, 3).Range.Text  = "869×3="
$t.Cell(1, 4).Range.Text  = "869×2="
$t.Cell(1, 5).Range.Text  = "764×5="

$t.Cell(5, 1).Range.Text  = "947×9="
$t.Cell(5, 2).Range.Text  = "558×3="
$t.Cell(5, 3).Range.Text  = "257×3="
$t.Cell(5, 4).Range.Text  = "316×6="
$t.Cell(5, 5).Range.Text  = "130×3="

$t.Cell(10, 1).Range.Text = "103×9="
$t.Cell(10, 2).Range.Text = "747×5="
$t.Cell(10, 3).Range.Text = "919×6="
$t.Cell(10, 4).Range.Text = "104×7="
$t.Cell(10, 5).Range.Text = "453×4="

$t.Cell(15, 1).Range.Text = "357×9="
$t.Cell(15, 2).Range.Text = "482×9="
$t.Cell(15, 3).Range.Text = "889×7="
$t.Cell(15, 4).Range.Text = "132×4="
$t.Cell(15, 5).Range.Text = "246×5="

$t.Cell(20, 1).Range.Text = "620×6="
$t.Cell(20, 2).Range.Text = "181×8="
$t.Cell(20, 3).Range.Text = "362×9="
$t.Cell(20, 4).Range.Text = "533×6="
$t.Cell(20, 5).Range.Text = "974×7="
